# SR [2022-08-21]: PS Class <- transfer of mgm logic started
#
# 1) Rename "Objects To Plan" -> "Team"
# 2) Populate the "Team" sheet with the team/property data table
# 3) Turn that range into a real Excel Table ("Table2")
# 4) Rename the "Tasks" table on "Planned Objects" -> "PlannedObjects"
# 5) Update selections / active sheet so "Team" becomes the active tab

$wb = $excel.ActiveWorkbook

$wsPlanned = $wb.Worksheets.Item(1)
$wsTeam = $wb.Worksheets.Item(2)

# --- Rename the second sheet -------------------------------------------------
$wsTeam.Name = "Team"

# --- Rename the existing table on the "Planned Objects" sheet ---------------
$plannedTable = $wsPlanned.ListObjects.Item(1)
$plannedTable.Name = "PlannedObjects"

# --- Header row ---------------------------------------------------------------
$wsTeam.Range("A3").Value = "Index"
$wsTeam.Range("B3").Value = "Property"
$wsTeam.Range("C3").Value = "Value"
$wsTeam.Range("D3").Value = "Email"
$wsTeam.Range("E3").Value = "StartDate"
$wsTeam.Range("F3").Value = "EndDate"
$wsTeam.Range("G3").Value = "Comment"

# --- Data rows ------------------------------------------------------------
$wsTeam.Range("A4").Value = 0
$wsTeam.Range("B4").Value = "Report version"
$wsTeam.Range("C4").Value = "'1.0"
$wsTeam.Range("D4").Value = "some_sample@email.com"
$wsTeam.Range("E4").Value = 44749
$wsTeam.Range("F4").Value = "'"
$wsTeam.Range("G4").Value = "'"

$wsTeam.Range("A5").Value = 1
$wsTeam.Range("B5").Value = "TL"
$wsTeam.Range("C5").Value = "Amade Wolfgang"
$wsTeam.Range("D5").Value = "some_sample@email.com"
$wsTeam.Range("E5").Value = 44749
$wsTeam.Range("F5").Value = "'"
$wsTeam.Range("G5").Value = "'"

$wsTeam.Range("A6").Value = 2
$wsTeam.Range("B6").Value = "BA"
$wsTeam.Range("C6").Value = "Ledowskykh Sergii"
$wsTeam.Range("D6").Value = "some_sample@email.com"
$wsTeam.Range("E6").Value = 44749
$wsTeam.Range("F6").Value = "'"
$wsTeam.Range("G6").Value = "'"

$wsTeam.Range("A7").Value = 3
$wsTeam.Range("B7").Value = "PO"
$wsTeam.Range("C7").Value = "Shelly Bengia"
$wsTeam.Range("D7").Value = "some_sample@email.com"
$wsTeam.Range("E7").Value = 44749
$wsTeam.Range("F7").Value = "'"
$wsTeam.Range("G7").Value = "'"

$wsTeam.Range("A8").Value = 4
$wsTeam.Range("B8").Value = "developer"
$wsTeam.Range("C8").Value = "Razumov Sergii"
$wsTeam.Range("D8").Value = "sergii_razumov@epam.com"
$wsTeam.Range("E8").Value = 44749
$wsTeam.Range("F8").Value = "'"
$wsTeam.Range("G8").Value = "'"

$wsTeam.Range("A9").Value = 5
$wsTeam.Range("B9").Value = "QA"
$wsTeam.Range("C9").Value = "Razumov Sergii"
$wsTeam.Range("D9").Value = "sergii_razumov@epam.com"
$wsTeam.Range("E9").Value = 44749
$wsTeam.Range("F9").Value = "'"
$wsTeam.Range("G9").Value = "'"

# --- Number formats --------------------------------------------------------
$wsTeam.Range("A4:A9").NumberFormat = "0"
$wsTeam.Range("E4:E9").NumberFormat = "dd\-mm\-yyyy"

# --- Column widths (matches the widths produced by Excel's "AutoFit") -----
$wsTeam.Range("A1").ColumnWidth = 7
$wsTeam.Range("B1").ColumnWidth = 12
$wsTeam.Range("C1").ColumnWidth = 15
$wsTeam.Range("D1").ColumnWidth = 23.166666666666668
$wsTeam.Range("E1").ColumnWidth = 10.333333333333334
$wsTeam.Range("F1").ColumnWidth = 10
$wsTeam.Range("G1").ColumnWidth = 10.666666666666666

# --- Turn the range into a Table --------------------------------------------
$teamTable = $wsTeam.ListObjects.Add(1, $wsTeam.Range("A3:G9"), [System.Type]::Missing, 1)
$teamTable.Name = "Table2"

# --- Selections / active sheet ----------------------------------------------
$wsPlanned.Range("E3").Select() | Out-Null
$wsTeam.Activate() | Out-Null
$wsTeam.Rows(1).Select() | Out-Null
